# Generate Report for Handback
# Renames the two e2e markdown source files tracked by the handback-status
# workbook (new GUID-based filenames) and refreshes the associated
# handoff/handback xliff filenames and timestamps produced by the newer run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Old / new identifiers
# ---------------------------------------------------------------------
$oldFile1 = "1e851055-11df-4ef1-924f-439760840548.md"
$newFile1 = "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.md"

$oldFile2 = "8c02821f-2870-4848-b5d5-3f80b9d200ad.md"
$newFile2 = "ffff24d56b07-c835-4088-9540-3f0a6c1c7f70.md"

$newXlfZhCn = "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.43a82b253913a1b3df6d97b1d24824a26ba7c7e5.zh-cn.xlf"
$newXlfDeDe = "26ee62c0-6b2f-4dd7-b704-2abfed8096a5.43a82b253913a1b3df6d97b1d24824a26ba7c7e5.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("G2").Value = "2016-09-04 15:06:11"
$wsOverview.Range("G3").Value = "2016-09-04 15:06:11"

# Hyperlinks on column B keep pointing at the original GitHub blob URLs
# (those targets are untouched) but their displayed text is refreshed to
# the new file names.
$ovHlRange = $wsOverview.Range("A1:G100")
$ovHlRange.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/1e851055-11df-4ef1-924f-439760840548.md", [Type]::Missing, [Type]::Missing, "e2e\$newFile1")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/8c02821f-2870-4848-b5d5-3f80b9d200ad.md", [Type]::Missing, [Type]::Missing, "e2e\$newFile2")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("G2").Value = $newXlfZhCn
$wsZhCn.Range("H2").Value = "2016-09-04 15:06:05"
$wsZhCn.Range("J2").Value = $newXlfZhCn
$wsZhCn.Range("K2").Value = "2016-09-04 15:06:37"

$wsZhCn.Range("G3").Value = $newXlfZhCn
$wsZhCn.Range("H3").Value = "2016-09-04 15:06:05"
$wsZhCn.Range("J3").Value = $newXlfZhCn
$wsZhCn.Range("K3").Value = "2016-09-04 15:06:37"

$zhHlRange = $wsZhCn.Range("A1:P100")
$zhHlRange.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/1e851055-11df-4ef1-924f-439760840548.md", [Type]::Missing, [Type]::Missing, $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/43fa034b8ec4d9a077192d366784f8fe5faa6005/e2e/1e851055-11df-4ef1-924f-439760840548.md", [Type]::Missing, [Type]::Missing, $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/8c02821f-2870-4848-b5d5-3f80b9d200ad.md", [Type]::Missing, [Type]::Missing, $newFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/43fa034b8ec4d9a077192d366784f8fe5faa6005/e2e/8c02821f-2870-4848-b5d5-3f80b9d200ad.md", [Type]::Missing, [Type]::Missing, $newFile2)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G2").Value = $newXlfDeDe
$wsDeDe.Range("H2").Value = "2016-09-04 15:06:11"
$wsDeDe.Range("J2").Value = $newXlfDeDe
$wsDeDe.Range("K2").Value = "2016-09-04 15:06:44"

$wsDeDe.Range("G3").Value = $newXlfDeDe
$wsDeDe.Range("H3").Value = "2016-09-04 15:06:11"
$wsDeDe.Range("J3").Value = $newXlfDeDe
$wsDeDe.Range("K3").Value = "2016-09-04 15:06:44"

$deHlRange = $wsDeDe.Range("A1:P100")
$deHlRange.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/1e851055-11df-4ef1-924f-439760840548.md", [Type]::Missing, [Type]::Missing, $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/95898753389455f1f52bca93fbe1ae2f5786b1a3/e2e/1e851055-11df-4ef1-924f-439760840548.md", [Type]::Missing, [Type]::Missing, $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54d5c65db63870f2fceeed6f0bcebdc693067ce6/e2e/8c02821f-2870-4848-b5d5-3f80b9d200ad.md", [Type]::Missing, [Type]::Missing, $newFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/95898753389455f1f52bca93fbe1ae2f5786b1a3/e2e/8c02821f-2870-4848-b5d5-3f80b9d200ad.md", [Type]::Missing, [Type]::Missing, $newFile2)

"done"
